$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dSF (column F) values for rows 2, 7, 8, 9 with repulled data
$ws.Range("F2").Value = 2
$ws.Range("F7").Value = 4
$ws.Range("F8").Value = -1
$ws.Range("F9").Value = -7
